# Tracker.xlsx update: python tasks added and tracker extended through 2018-08-16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 17: it previously held only A17/B17 (with the typo "pyhton").
#     Correct B17 to "python" and add the missing C17 description so the
#     row becomes a full 3-column record like the others.
$ws.Range("B17").Value = "python"
$ws.Range("C17").Value = "Hankerank tasks on python"

# --- Append two more tracker rows (18 and 19), copying the formatting
#     from row 17 so number formats / fills / borders stay consistent.
$ws.Range("A17:C17").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)
$ws.Range("A17:C17").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)

$ws.Range("A18").Value = 43327
$ws.Range("B18").Value = "python"
$ws.Range("C18").Value = "Hankerank tasks on python"

$ws.Range("A19").Value = 43328
$ws.Range("B19").Value = "python"
$ws.Range("C19").Value = "Hankerank tasks on python"

# --- Update the active selection to the last entry, matching the tracker's
#     "current position" convention used throughout this workbook.
$ws.Range("C19").Select()
